$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 350.4
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 576
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 576
$ws.Range("M12").Value = -30
$ws.Range("N12").Value = -916
$ws.Range("H17").Value = 1254.16
$ws.Range("J17").Value = 1254.16
$ws.Range("L17").Value = 3762.48
$ws.Range("N17").Value = -4098.48
$ws.Range("H29").Value = 485
$ws.Range("I29").Value = 146.66667
$ws.Range("K29").Value = 440.00001
$ws.Range("M29").Value = -159.00001
$ws.Range("H92").Value = 443.33334
$ws.Range("I92").Value = 223.33333
$ws.Range("J92").Value = 663.3333
$ws.Range("K92").Value = 223.33333
$ws.Range("L92").Value = 663.3333
$ws.Range("M92").Value = 1024.66667
$ws.Range("N92").Value = -3159.3333
$ws.Range("H108").Value = 40684
$ws.Range("J108").Value = 40684
$ws.Range("L108").Value = 40684
$ws.Range("N108").Value = -48364
$ws.Range("H129").Value = 995.2941
$ws.Range("I129").Value = 677.9
$ws.Range("J129").Value = 1127.5416
$ws.Range("K129").Value = 2033.7
$ws.Range("L129").Value = 3382.6248
$ws.Range("M129").Value = 2966.3
$ws.Range("N129").Value = -13382.6248
$ws.Range("H132").Value = 1916.9744
$ws.Range("I132").Value = 1601.3334
$ws.Range("J132").Value = 3653
$ws.Range("K132").Value = 4804.0002
$ws.Range("L132").Value = 10959
$ws.Range("M132").Value = -2274.0002
$ws.Range("N132").Value = -16019
$ws.Range("H138").Value = 1886.27
$ws.Range("I138").Value = 1380.9333
$ws.Range("J138").Value = 2102.8428
$ws.Range("K138").Value = 4142.7999
$ws.Range("L138").Value = 6308.528399999999
$ws.Range("M138").Value = 997.2001
$ws.Range("N138").Value = -16588.5284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4166.6665
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 311998
$ws.Range("I26").Value = 311998
$ws.Range("K26").Value = 311998
$ws.Range("M26").Value = -311706
$ws.Range("H96").Value = 144928.38
$ws.Range("I96").Value = 181571.17
$ws.Range("K96").Value = 181571.17
$ws.Range("M96").Value = -178825.17
$ws.Range("H102").Value = 500000
$ws.Range("I102").Value = 500000
$ws.Range("K102").Value = 500000
$ws.Range("M102").Value = -496755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21279314
$ws.Range("I31").Value = 45456804
$ws.Range("J31").Value = 3121.88
$ws.Range("K31").Value = 45456804
$ws.Range("L31").Value = 3121.88
$ws.Range("M31").Value = -45456509
$ws.Range("N31").Value = -3711.88
$ws.Range("H34").Value = 21279314
$ws.Range("I34").Value = 45456804
$ws.Range("J34").Value = 3121.88
$ws.Range("K34").Value = 45456804
$ws.Range("L34").Value = 3121.88
$ws.Range("M34").Value = -45456602
$ws.Range("N34").Value = -3525.88
$ws.Range("H86").Value = 5343.6
$ws.Range("I86").Value = 5859
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 5859
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -4736
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 5343.6
$ws.Range("I89").Value = 5859
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 29295
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -23679
$ws.Range("N89").Value = -36232
$ws.Range("H97").Value = 40189.668
$ws.Range("J97").Value = 40189.668
$ws.Range("L97").Value = 40189.668
$ws.Range("N97").Value = -42171.668
$ws.Range("H132").Value = 451896.97
$ws.Range("I132").Value = 615418.7
$ws.Range("J132").Value = 2212.25
$ws.Range("K132").Value = 1846256.1
$ws.Range("L132").Value = 6636.75
$ws.Range("M132").Value = -1843726.1
$ws.Range("N132").Value = -11696.75
$ws.Range("H134").Value = 1641.1177
$ws.Range("I134").Value = 1573.75
$ws.Range("K134").Value = 4721.25
$ws.Range("M134").Value = -2186.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 840210.5
$ws.Range("I12").Value = 110.25
$ws.Range("J12").Value = 1017073.75
$ws.Range("K12").Value = 330.75
$ws.Range("L12").Value = 3051221.25
$ws.Range("M12").Value = -157.75
$ws.Range("N12").Value = -3051567.25
$ws.Range("H33").Value = 1108.5
$ws.Range("J33").Value = 1557.4286
$ws.Range("L33").Value = 9344.571599999999
$ws.Range("N33").Value = -9910.571599999999
$ws.Range("H87").Value = 3166.6667
$ws.Range("I87").Value = 1800
$ws.Range("K87").Value = 5400
$ws.Range("M87").Value = -4152
$ws.Range("H88").Value = 140000
$ws.Range("J88").Value = 140000
$ws.Range("L88").Value = 420000
$ws.Range("N88").Value = -420856
$ws.Range("H90").Value = 3166.6667
$ws.Range("I90").Value = 1800
$ws.Range("K90").Value = 16200
$ws.Range("M90").Value = -9960
$ws.Range("H91").Value = 140000
$ws.Range("J91").Value = 140000
$ws.Range("L91").Value = 420000
$ws.Range("N91").Value = -422964
$ws.Range("H114").Value = 400.95456
$ws.Range("J114").Value = 426.64285
$ws.Range("L114").Value = 1279.92855
$ws.Range("N114").Value = -7787.928550000001
$ws.Range("H131").Value = 3479.3462
$ws.Range("I131").Value = 12896
$ws.Range("J131").Value = 1767.2273
$ws.Range("K131").Value = 38688
$ws.Range("L131").Value = 5301.6819
$ws.Range("M131").Value = -33648
$ws.Range("N131").Value = -15381.6819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 15501.143
$ws.Range("I97").Value = 19037.143
$ws.Range("J97").Value = 1357.1428
$ws.Range("K97").Value = 19037.143
$ws.Range("L97").Value = 1357.1428
$ws.Range("M97").Value = -18541.143
$ws.Range("N97").Value = -2349.1428
$ws.Range("H122").Value = 3110.7222
$ws.Range("I122").Value = 3400.9285
$ws.Range("J122").Value = 2095
$ws.Range("K122").Value = 10202.7855
$ws.Range("L122").Value = 6285
$ws.Range("M122").Value = -7752.7855
$ws.Range("N122").Value = -11185
$ws.Range("H123").Value = 13036.4375
$ws.Range("J123").Value = 13036.4375
$ws.Range("L123").Value = 13036.4375
$ws.Range("N123").Value = -17936.4375
$ws.Range("H132").Value = 2449.5652
$ws.Range("I132").Value = 1622.3334
$ws.Range("J132").Value = 4000.625
$ws.Range("K132").Value = 4867.0002
$ws.Range("L132").Value = 12001.875
$ws.Range("M132").Value = -2337.0002
$ws.Range("N132").Value = -17061.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 32250
$ws.Range("I5").Value = 50000
$ws.Range("J5").Value = 14500
$ws.Range("K5").Value = 50000
$ws.Range("L5").Value = 14500
$ws.Range("M5").Value = -49887
$ws.Range("N5").Value = -14726
$ws.Range("H93").Value = 1599.8
$ws.Range("I93").Value = 1250
$ws.Range("K93").Value = 1250
$ws.Range("M93").Value = -2
$ws.Range("H136").Value = 2763.577
$ws.Range("I136").Value = 2410.5417
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 7231.625100000001
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -4681.625100000001
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 27500
$ws.Range("H124").Value = 66000
$ws.Range("J124").Value = 66000
$ws.Range("L124").Value = 66000
$ws.Range("N124").Value = -75820

